$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove stray comma from the slime-swamp description (row 5, column G)
$ws.Range("G5").Value2 = "슬라임들이 살아가는 늪지. 다만 늪은 없다. 예전에는 있었다고 한다."

# Remove stray comma from the wolf-area description (row 9, column G)
$ws.Range("G9").Value2 = "늑대들의 주신 암월랑의 출몰지역. 별밤의 황홀함에 현혹되지 말어라. 그건 네 앞에 있다."

# Move the active selection to G12
$ws.Range("G12").Select() | Out-Null
